# storyTrigger.xlsx edit: add a new "npcNotInTeam" column (string) right
# after "npcInTeam" (old column N), shifting repeatable/forbiddenStoryId/
# priority/money one column to the right, and append 3 new template rows
# (5,6,7) copied from the row-3 pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column at N (14), shifting N..Q -> O..R -------------
$ws.Columns.Item(14).Insert()

# --- 2. New header cell N1: "npcNotInTeam" with a bold accent font -------
$hdr = $ws.Cells.Item(1, 14)
$hdr.Value = "npcNotInTeam"
$hdr.Font.Bold = $true
$hdr.Font.Name = "Helvetica Neue"
$hdr.Font.Color = 3355443   # RGB(51,51,51) = FF333333

$ws.Rows.Item(1).RowHeight = 16

# --- 3. New "type" marker cell N2: "string" -------------------------------
$ws.Cells.Item(2, 14).Value = "string"

# --- 4. Data rows 3-4 for the new column: all zero ------------------------
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(4, 14).Value = 0

# --- 5. Move the "repeatable" header comment from its old cell (now N1,
#        empty) to its new location O1, preserving the text/author ---------
$oldCommentCell = $ws.Cells.Item(1, 14)
if ($oldCommentCell.Comment) {
    $commentText = $oldCommentCell.Comment.Text()
    $oldCommentCell.Comment.Delete()
    $ws.Cells.Item(1, 15).AddComment($commentText)
}

# --- 6. Append three new template rows (5,6,7), matching row 3's pattern,
#        with F6 = 2 per the source data ------------------------------------
$newRows = @(
    @{ Row = 5; A = 3; B = 1; F = 0 },
    @{ Row = 6; A = 4; B = 1; F = 2 },
    @{ Row = 7; A = 5; B = 1; F = 0 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A    # A storyId
    $ws.Cells.Item($row, 2).Value = $r.B    # B locked
    $ws.Cells.Item($row, 3).Value = 2       # C heroId
    $ws.Cells.Item($row, 4).Value = 0       # D cityId
    $ws.Cells.Item($row, 5).Value = 0       # E prefixStoryId
    $ws.Cells.Item($row, 6).Value = $r.F    # F buildingId
    $ws.Cells.Item($row, 7).Value = 0       # G month
    $ws.Cells.Item($row, 8).Value = 0       # H day
    $ws.Cells.Item($row, 9).Value = 0       # I year
    $ws.Cells.Item($row, 10).Value = 0      # J cityPercentage
    $ws.Cells.Item($row, 11).Value = 0      # K commerce
    $ws.Cells.Item($row, 12).Value = 0      # L military
    $ws.Cells.Item($row, 13).Value = 0      # M npcInTeam
    $ws.Cells.Item($row, 14).Value = 0      # N npcNotInTeam
    $ws.Cells.Item($row, 15).Value = 0      # O repeatable
    $ws.Cells.Item($row, 16).Value = 0      # P forbiddenStoryId
    $ws.Cells.Item($row, 17).Value = 1000   # Q priority
    $ws.Cells.Item($row, 18).Value = 0      # R money
}

# --- 7. Misc cosmetic bits that mirror the saved-from-Excel state ---------
$ws.PageSetup.Orientation = 1   # xlPortrait

$ws.Range("L14").Select() | Out-Null
